# Updated Argent (Silver) prices: append the new daily price row (2025-05-30)
# to the bottom of the "Prices" sheet, carrying forward the prior day's
# values (as in the source diff) for columns B:J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 90

# The sheet stores every value as text (dates, numbers, thousands-separated
# strings like "5,356" alike). Force the new cells to Text format first so
# Excel's automatic type detection doesn't turn "2025-05-30" into a date
# serial or "35.5" into a real number, then set the literal values.
$rng = $ws.Range("A" + $newRow + ":J" + $newRow)
$rng.NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-05-30"
$ws.Range("B" + $newRow).Value = "35.5"
$ws.Range("C" + $newRow).Value = "35.01"
$ws.Range("D" + $newRow).Value = "0.94"
$ws.Range("E" + $newRow).Value = "0.253"
$ws.Range("F" + $newRow).Value = "0.09"
$ws.Range("G" + $newRow).Value = "5,356"
$ws.Range("H" + $newRow).Value = "8,019"
$ws.Range("I" + $newRow).Value = "8,069"
$ws.Range("J" + $newRow).Value = "7.2186"

# Restore the plain/default cell style so the new row doesn't carry a
# leftover "Text" number-format style (matches the rest of the sheet, which
# uses the default style for every cell).
$rng.Style = "Normal"
